$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append two new test-run rows at the bottom of the history
# table (rows 60 and 61).
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")
$amsinDateFmt = $wsAmsin.Range("B59").NumberFormat

$wsAmsin.Range("A60").Value = "'2023-02-17"
$wsAmsin.Range("B60").NumberFormat = $amsinDateFmt
$wsAmsin.Range("B60").Value = 44974.43709384259
$wsAmsin.Range("C60").Value = "173cyclefst"
$wsAmsin.Range("D60").Value = 124
$wsAmsin.Range("E60").Value = 123
$wsAmsin.Range("F60").Value = 1
$wsAmsin.Range("G60").Value = 2.05

$wsAmsin.Range("A61").Value = "'2023-02-20"
$wsAmsin.Range("B61").NumberFormat = $amsinDateFmt
$wsAmsin.Range("B61").Value = 44977.40212696759
$wsAmsin.Range("C61").Value = "173fnlrun"
$wsAmsin.Range("D61").Value = 124
$wsAmsin.Range("E61").Value = 113
$wsAmsin.Range("F61").Value = 11
$wsAmsin.Range("G61").Value = 4.33

# ---------------------------------------------------------------------------
# Sheet "BETA": append one new test-run row at the bottom of the history
# table (row 32).
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
$betaDateFmt = $wsBeta.Range("B31").NumberFormat

$wsBeta.Range("A32").Value = "'2023-02-20"
$wsBeta.Range("B32").NumberFormat = $betaDateFmt
$wsBeta.Range("B32").Value = 44977.58745057871
$wsBeta.Range("C32").Value = "173beta"
$wsBeta.Range("D32").Value = 124
$wsBeta.Range("E32").Value = 123
$wsBeta.Range("F32").Value = 1
$wsBeta.Range("G32").Value = 1.66

# ---------------------------------------------------------------------------
# Sheet "AMS": row 35 was missing the standard cell formatting the rest of
# the table uses, and its run-time value had drifted slightly - bring it
# back in line with the rest of the sheet, then append the newest run as
# row 36.
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")
$amsDateFmt = $wsAms.Range("B34").NumberFormat

$wsAms.Range("A35").NumberFormat = "General"
$wsAms.Range("C35").NumberFormat = "General"
$wsAms.Range("D35").NumberFormat = "General"
$wsAms.Range("E35").NumberFormat = "General"
$wsAms.Range("F35").NumberFormat = "General"
$wsAms.Range("G35").NumberFormat = "General"
$wsAms.Range("B35").Value = 44946.89168506944

$wsAms.Range("A36").Value = "'2023-02-20"
$wsAms.Range("B36").NumberFormat = $amsDateFmt
$wsAms.Range("B36").Value = 44977.83446887095
$wsAms.Range("C36").Value = "live173"
$wsAms.Range("D36").Value = 124
$wsAms.Range("E36").Value = 123
$wsAms.Range("F36").Value = 1
$wsAms.Range("G36").Value = 1.78
